$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.997390666666667
$ws.Range("H2").Value = 8.992172
$ws.Range("I2").Value = 0.01499229953737403
$ws.Range("J2").Value = 0.01499229953737403
$ws.Range("M2").Value = 0.428743
$ws.Range("N2").Value = 1.286229
$ws.Range("O2").Value = 0.00412050394863168
$ws.Range("P2").Value = 0.00412050394863168
$ws.Range("Q2").Value = 1.285110266598667
$ws.Range("R2").Value = 11.565992399388
$ws.Range("S2").Value = 0.0000617758294428186
$ws.Range("T2").Value = 0.0000617758294428186

$ws.Range("G3").Value = 2.997390666666667
$ws.Range("H3").Value = 8.992172
$ws.Range("I3").Value = 0.01499229953737403
$ws.Range("J3").Value = 0.01499229953737403
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("N3").Value = 240.678711
$ws.Range("O3").Value = 0.7710272268990069
$ws.Range("P3").Value = 0.7710272268990069
$ws.Range("Q3").Value = 240.469374005588
$ws.Range("R3").Value = 2164.224366050292
$ws.Range("S3").Value = 0.01155947113714076
$ws.Range("T3").Value = 0.01155947113714076

$ws.Range("G4").Value = 2.997390666666667
$ws.Range("H4").Value = 8.992172
$ws.Range("I4").Value = 0.01499229953737403
$ws.Range("J4").Value = 0.01499229953737403
$ws.Range("M4").Value = 23.39612766666667
$ws.Range("N4").Value = 70.188383
$ws.Range("O4").Value = 0.2248522691523614
$ws.Range("P4").Value = 0.2248522691523614
$ws.Range("Q4").Value = 70.12733470420845
$ws.Range("R4").Value = 631.1460123378761
$ws.Range("S4").Value = 0.003371052570790449
$ws.Range("T4").Value = 0.003371052570790449

$ws.Range("I5").Value = 0.825470460014473
$ws.Range("J5").Value = 0.825470460014473
$ws.Range("M5").Value = 0.428743
$ws.Range("N5").Value = 1.286229
$ws.Range("O5").Value = 0.00412050394863168
$ws.Range("P5").Value = 0.00412050394863168
$ws.Range("Q5").Value = 70.75769532846002
$ws.Range("R5").Value = 636.8192579561401
$ws.Range("S5").Value = 0.003401354289968446
$ws.Range("T5").Value = 0.003401354289968446

$ws.Range("I6").Value = 0.825470460014473
$ws.Range("J6").Value = 0.825470460014473
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("N6").Value = 240.678711
$ws.Range("O6").Value = 0.7710272268990069
$ws.Range("P6").Value = 0.7710272268990069
$ws.Range("S6").Value = 0.6364601996720067
$ws.Range("T6").Value = 0.6364601996720067

$ws.Range("I7").Value = 0.825470460014473
$ws.Range("J7").Value = 0.825470460014473
$ws.Range("M7").Value = 23.39612766666667
$ws.Range("N7").Value = 70.188383
$ws.Range("O7").Value = 0.2248522691523614
$ws.Range("P7").Value = 0.2248522691523614
$ws.Range("Q7").Value = 3861.185076616421
$ws.Range("R7").Value = 34750.66568954779
$ws.Range("S7").Value = 0.1856089060524979
$ws.Range("T7").Value = 0.1856089060524979

$ws.Range("G8").Value = 31.89607
$ws.Range("H8").Value = 95.68821
$ws.Range("I8").Value = 0.159537240448153
$ws.Range("J8").Value = 0.159537240448153
$ws.Range("M8").Value = 0.428743
$ws.Range("N8").Value = 1.286229
$ws.Range("O8").Value = 0.00412050394863168
$ws.Range("P8").Value = 0.00412050394863168
$ws.Range("Q8").Value = 13.67521674001
$ws.Range("R8").Value = 123.07695066009
$ws.Range("S8").Value = 0.0006573738292204162
$ws.Range("T8").Value = 0.0006573738292204162

$ws.Range("G9").Value = 31.89607
$ws.Range("H9").Value = 95.68821
$ws.Range("I9").Value = 0.159537240448153
$ws.Range("J9").Value = 0.159537240448153
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("N9").Value = 240.678711
$ws.Range("O9").Value = 0.7710272268990069
$ws.Range("P9").Value = 0.7710272268990069
$ws.Range("Q9").Value = 2558.901671188589
$ws.Range("R9").Value = 23030.11504069731
$ws.Range("S9").Value = 0.1230075560898595
$ws.Range("T9").Value = 0.1230075560898595

$ws.Range("G10").Value = 31.89607
$ws.Range("H10").Value = 95.68821
$ws.Range("I10").Value = 0.159537240448153
$ws.Range("J10").Value = 0.159537240448153
$ws.Range("M10").Value = 23.39612766666667
$ws.Range("N10").Value = 70.188383
$ws.Range("O10").Value = 0.2248522691523614
$ws.Range("P10").Value = 0.2248522691523614
$ws.Range("Q10").Value = 746.2445257849366
$ws.Range("R10").Value = 6716.20073206443
$ws.Range("S10").Value = 0.0358723105290731
$ws.Range("T10").Value = 0.0358723105290731
